$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -6544.271299816632
$ws.Range("C2").Value = 13709.61546051708
$ws.Range("D2").Value = -6544.271301772727
$ws.Range("E2").Value = -621.0728589275514
$ws.Range("F2").Value = 56.87534216409858
$ws.Range("G2").Value = 45.34803821248835
$ws.Range("H2").Value = 57.19007467166983
$ws.Range("I2").Value = 45.11578120160704
$ws.Range("J2").Value = 56.87534216421642
$ws.Range("K2").Value = 45.34803821349107
$ws.Range("L2").Value = 54.34895899973804
$ws.Range("O2").Value = 45.41483505640417
$ws.Range("P2").Value = 54.34895900253247
$ws.Range("R2").Value = 8.1102736936438
$ws.Range("S2").Value = -16.22054739033433
$ws.Range("T2").Value = 8.11027369669053
$ws.Range("X2").Value = -134.2333058426104
$ws.Range("Y2").Value = -175.2596612421746
$ws.Range("Z2").Value = -134.2333058426103
$ws.Range("AE2").Value = -8.1102736936438
$ws.Range("AF2").Value = 8.11027369669053
$ws.Range("AG2").Value = 8.1102736936438
$ws.Range("AH2").Value = -16.22054739033433
$ws.Range("AI2").Value = 8.11027369669053
$ws.Range("AJ2").Value = 8.1102736936438
$ws.Range("AK2").Value = -8.11027369669053
$ws.Range("AL2").Value = 20.51317769978211
$ws.Range("AM2").Value = -20.5131776997821
$ws.Range("AN2").Value = -134.2333058426104
$ws.Range("AO2").Value = -175.2596612421746
$ws.Range("AP2").Value = -134.2333058426103
$ws.Range("AQ2").Value = -20.51317769978211
$ws.Range("AR2").Value = 20.5131776997821
$ws.Range("AS2").Value = 56.87534216409858
$ws.Range("AT2").Value = 56.87534216409858
$ws.Range("AU2").Value = 57.19007467166989
$ws.Range("AV2").Value = 57.19007467166989
$ws.Range("AW2").Value = 57.19007467166983
$ws.Range("AX2").Value = 56.87534216421642
$ws.Range("AY2").Value = 56.87534216421642
$ws.Range("AZ2").Value = 45.34803821248835
$ws.Range("BA2").Value = 45.34803821248835
$ws.Range("BB2").Value = 45.11578120160704
$ws.Range("BC2").Value = 45.11578120106589
$ws.Range("BD2").Value = 45.11578120214858
$ws.Range("BE2").Value = 45.34803821349107
$ws.Range("BF2").Value = 45.34803821349112
$ws.Range("BG2").Value = 54.34895899973804
$ws.Range("BJ2").Value = 45.41483505640417
$ws.Range("BK2").Value = 54.34895900253247
